# Replicate the workbook edit described in the commit:
#  - Add a new worksheet named "Sheet1" after the existing "Sheet 1"
#  - Copy the top of the ranking table (header + first 5 rows, A1:D6)
#    from "Sheet 1" into the new sheet, preserving original text typing
#  - New sheet becomes the active / selected tab, with E12 selected
#  - "Sheet 1" keeps A1:D6 selected (no longer the active tab)
#  - Tidy up row heights on "Sheet 1" (drop per-row overrides)

$wb = $excel.ActiveWorkbook
$sheet1 = $wb.Worksheets.Item(1)

# Select A1:D6 on the original sheet while it is still the active sheet,
# matching the final selection state recorded for "Sheet 1".
$sheet1.Range("A1:D6").Select()

# Drop the per-row height overrides on "Sheet 1" (now a uniform height).
$sheet1.Rows("1:51").AutoFit()

# Insert the new worksheet right after "Sheet 1" and rename it.
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $sheet1)
$newSheet.Name = "Sheet1"

# Copy A1:D6 (header + first five ranked companies) across, preserving the
# original cell typing (numbers-stored-as-text) via Copy/PasteSpecial
# instead of a Value assignment (which would coerce numeric-looking text).
$sheet1.Range("A1:D6").Copy()
$newSheet.Range("A1").PasteSpecial()

# Final selection on the new sheet.
$newSheet.Range("E12").Select()
